$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
# B2 holds a purely-numeric-looking string ("4"); force Text storage via
# NumberFormat, assign, then strip the format again so the cell ends up
# back at the default style but keeps its literal-text value/type.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "fdhxfgh"
$ws.Range("D2").Value = "vcgh"

# --- Row 3 updates ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "6"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Value = "xfgh"
$ws.Range("D3").Value = "fgh"

# --- New row 4 ---
# Clone A3's full formatting (border/font/alignment) onto A4 via Copy,
# then overwrite the copied value with the real one.
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 3

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "8"
$ws.Range("B4").ClearFormats()

$ws.Range("C4").Value = "fgjn"
$ws.Range("D4").Value = "fgjn"
